$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting the existing rows 4..19 down to 5..20.
# (The new row inherits the formatting, including the date style, from the row above.)
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with this week's data.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44487
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112013
$ws.Range("G4").Value = "Alcachofa"
$ws.Range("H4").Value = "Madrigal"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 11500
$ws.Range("N4").Value = "$/caja 40 unidades"
$ws.Range("O4").Value = "Provincia del Elquí"
$ws.Range("P4").Value = 288
$ws.Range("Q4").Value = 40
$ws.Range("R4").Value = "Hortaliza"
